# Reporte-Prueba borrador.xlsx edit
# - Adds 36 new nutrient columns (N:AW) to the header row
# - Renames the encuestador "Donato" -> "Leandro Donato" in column A (rows 2-5)
# - Updates the nutrient sample values for rows 3, 4 and 5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns N1:AW1 ------------------------------------------------
$newHeaders = @(
    "Colesterol",
    "Fibra Alimentaria",
    "Sodio",
    "Agua",
    "Vitamina A",
    "Vitamina B6",
    "Vitamina B12",
    "Vitamina C",
    "Vitamina D",
    "Vitamina E",
    "Vitamina K",
    "Almidón",
    "Lactosa",
    "Alcohol",
    "Cafeína",
    "Azúcares",
    "Calcio",
    "Hierro",
    "Magnesio",
    "Fósforo",
    "Cinc",
    "Cobre",
    "Fluor",
    "Manganeso",
    "Selenio",
    "Tiamina",
    "Ácido Pantetónico",
    "Riboflavina",
    "Niacina",
    "Folato",
    "Ácido Fólico",
    "Gasas Trans",
    "Grasas Monoinsaturadas",
    "Grasas Poliinsaturadas",
    "Cloruro",
    "Caroteno"
)

$startCol = 14 # column N
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $newHeaders[$i]
}

# --- Encuestador name update (A2:A5) ------------------------------------------
$ws.Range("A2").Value = "Leandro Donato"
$ws.Range("A3").Value = "Leandro Donato"
$ws.Range("A4").Value = "Leandro Donato"
$ws.Range("A5").Value = "Leandro Donato"

# --- Row 3 nutrient values ------------------------------------------------------
$ws.Range("H3").Value = 150
$ws.Range("J3").Value = 30.0
$ws.Range("K3").Value = 45.0
$ws.Range("L3").Value = 60.0
$ws.Range("M3").Value = 75.0

# --- Row 4 nutrient values ------------------------------------------------------
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 20
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = 40

# --- Row 5 nutrient values ------------------------------------------------------
$ws.Range("H5").Value = 10
$ws.Range("J5").Value = 2.0
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 4.0
$ws.Range("M5").Value = 5.0
